$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.373617649078369
$ws.Range("B1").Value = 2.042556285858154
$ws.Range("C1").Value = 3.915194511413574
$ws.Range("D1").Value = 1.077288866043091
$ws.Range("E1").Value = 0.7182868123054504
